$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename IMAGE1 value and add a new IMAGE2 value
$ws.Range("J2").Value = "perceuse-01.jpg"
$ws.Range("J3").Value = "perceuse-01.jpg"
$ws.Range("K2").Value = "perceuse-02.jpg"
$ws.Range("K3").Value = "perceuse-02.jpg"

# Match Excel's bestFit column width recalculation for the updated columns
# (IMAGE1/IMAGE2 now hold longer "perceuse-0x.jpg" filenames)
$ws.Columns.Item(10).ColumnWidth = 11.2
$ws.Columns.Item(11).ColumnWidth = 11.1

# Update selection
$ws.Range("J8").Select() | Out-Null
